$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet and update the "as of" label (new data through 2021-09-19)
$ws.Name = "Through 2021-09-19"
$ws.Range("B1").Value = "September 2021 (through September 19)"

# Garfield Park (row 2)
$ws.Range("B2").Value = 13
$ws.Range("K2").Value = 5
$ws.Range("T2").Value = 3

# North Lawndale (row 3)
$ws.Range("B3").Value = 6
$ws.Range("BD3").Value = 2

# Austin (row 5)
$ws.Range("K5").Value = 8
$ws.Range("T5").Value = 2
$ws.Range("AC5").Value = 6

# Roseland (row 6)
$ws.Range("B6").Value = 5

# Little Italy, UIC (row 11)
$ws.Range("K11").Value = 2

# South Shore (row 17)
$ws.Range("B17").Value = 2
$ws.Range("K17").Value = 3

# Wicker Park (row 19)
$ws.Range("B19").Value = 4

# Englewood (row 20)
$ws.Range("B20").Value = 2

# United Center (row 23)
$ws.Range("T23").Value = 1

# Logan Square (row 28)
$ws.Range("AC28").Value = 2

# Belmont Cragin (row 29)
$ws.Range("B29").Value = 1
$ws.Range("AU29").Value = 1

# Lake View (row 33)
$ws.Range("B33").Value = 3

# Hyde Park (row 34)
$ws.Range("B34").Value = 1

# Irving Park (row 53)
$ws.Range("B53").Value = 1

# Grand Crossing (row 55)
$ws.Range("K55").Value = 6
$ws.Range("AL55").Value = 1

# Edgewater (row 56)
$ws.Range("B56").Value = 2

# Pullman (row 91)
$ws.Range("B91").Value = 1

# Sheffield & DePaul (row 96)
$ws.Range("K96").Value = 1
